$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column K: "carrera_id"
$ws.Range("K1").Value = "carrera_id"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Data rows for new column K
$ws.Range("K2").Value = 1
$ws.Range("A2").Copy()
$ws.Range("K2").PasteSpecial(-4122)

$ws.Range("K3").Value = 2
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# Update selection to match target (K4)
$ws.Range("K4").Select()
